$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Current (DC Units)" values in column J for rows 8 and 9
$ws.Range("J8").Value = 6
$ws.Range("J9").Value = 6

# Update the active selection to match the saved view state
$ws.Range("I9").Select()
